# Update the AgTests (F) and AgPosit (G) rolling/cumulative figures
# for recent days, and append the new day (row 735) of data,
# matching the commit "Updated: st 10. 03. 2022".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(635, 6).Value = 83641
$ws.Cells.Item(643, 6).Value = 43479
$ws.Cells.Item(644, 6).Value = 36899
$ws.Cells.Item(651, 6).Value = 37139
$ws.Cells.Item(656, 6).Value = 52567
$ws.Cells.Item(663, 6).Value = 37265
$ws.Cells.Item(670, 6).Value = 52651
$ws.Cells.Item(677, 6).Value = 56201
$ws.Cells.Item(684, 6).Value = 57278
$ws.Cells.Item(691, 6).Value = 62645
$ws.Cells.Item(695, 6).Value = 37313
$ws.Cells.Item(697, 6).Value = 29115
$ws.Cells.Item(697, 7).Value = 3062
$ws.Cells.Item(698, 6).Value = 70847
$ws.Cells.Item(698, 7).Value = 5828
$ws.Cells.Item(700, 6).Value = 43780
$ws.Cells.Item(700, 7).Value = 4327
$ws.Cells.Item(704, 6).Value = 25097
$ws.Cells.Item(704, 7).Value = 3715
$ws.Cells.Item(705, 6).Value = 56016
$ws.Cells.Item(705, 7).Value = 6311
$ws.Cells.Item(707, 6).Value = 38823
$ws.Cells.Item(711, 6).Value = 22592
$ws.Cells.Item(711, 7).Value = 3819
$ws.Cells.Item(712, 6).Value = 51369
$ws.Cells.Item(712, 7).Value = 6317
$ws.Cells.Item(713, 6).Value = 37314
$ws.Cells.Item(713, 7).Value = 4758
$ws.Cells.Item(714, 6).Value = 32516
$ws.Cells.Item(714, 7).Value = 3995
$ws.Cells.Item(715, 6).Value = 31834
$ws.Cells.Item(715, 7).Value = 3569
$ws.Cells.Item(716, 6).Value = 29762
$ws.Cells.Item(716, 7).Value = 3683
$ws.Cells.Item(717, 6).Value = 12563
$ws.Cells.Item(717, 7).Value = 2135
$ws.Cells.Item(718, 6).Value = 17108
$ws.Cells.Item(718, 7).Value = 2865
$ws.Cells.Item(719, 6).Value = 44629
$ws.Cells.Item(719, 7).Value = 5225
$ws.Cells.Item(720, 6).Value = 31238
$ws.Cells.Item(720, 7).Value = 3517
$ws.Cells.Item(721, 6).Value = 27923
$ws.Cells.Item(721, 7).Value = 3147
$ws.Cells.Item(722, 6).Value = 27942
$ws.Cells.Item(722, 7).Value = 2881
$ws.Cells.Item(723, 6).Value = 22519
$ws.Cells.Item(723, 7).Value = 2742
$ws.Cells.Item(724, 6).Value = 9373
$ws.Cells.Item(724, 7).Value = 1506
$ws.Cells.Item(725, 6).Value = 12719
$ws.Cells.Item(725, 7).Value = 2079
$ws.Cells.Item(726, 6).Value = 35772
$ws.Cells.Item(726, 7).Value = 4100
$ws.Cells.Item(727, 6).Value = 24946
$ws.Cells.Item(727, 7).Value = 2777
$ws.Cells.Item(728, 6).Value = 24532
$ws.Cells.Item(728, 7).Value = 2584
$ws.Cells.Item(729, 6).Value = 22952
$ws.Cells.Item(729, 7).Value = 2467
$ws.Cells.Item(730, 6).Value = 19119
$ws.Cells.Item(730, 7).Value = 2261
$ws.Cells.Item(731, 6).Value = 8525
$ws.Cells.Item(731, 7).Value = 1297
$ws.Cells.Item(732, 6).Value = 11614
$ws.Cells.Item(732, 7).Value = 1859
$ws.Cells.Item(733, 6).Value = 30761
$ws.Cells.Item(733, 7).Value = 3572
$ws.Cells.Item(734, 6).Value = 22282
$ws.Cells.Item(734, 7).Value = 2441

# Append the new row of data (row 735) for date 2022-03-09 (serial 44629)
$ws.Cells.Item(735, 1).Value = 44629
$ws.Cells.Item(735, 1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(735, 2).Value = 1546510
$ws.Cells.Item(735, 3).Value = 20595
$ws.Cells.Item(735, 4).Value = 12066
$ws.Cells.Item(735, 5).Value = 18817
$ws.Cells.Item(735, 6).Value = 12272
$ws.Cells.Item(735, 7).Value = 1560

